$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old header row (row 8: Index/Date/Value) and shift data up.
# First, delete the last data row (row 18) since the new layout only needs 17 rows total.
$ws.Rows("18").Delete()

# Row 8 becomes the first data row (Index=1, Date=31May2020 2300, Value=-1)
$ws.Range("A8").Value = 1
$ws.Range("B8").Value = "31May2020  2300"
$ws.Range("C8").Value = -1

# Row 9: Index=2, Date=31May2020 2315, Value = C8+1 (first real formula)
$ws.Range("A9").Value = 2
$ws.Range("B9").Value = "31May2020  2315"
$ws.Range("C9").Formula = "=C8+1"

# Rows 10-17: Index 3..10, Dates - set all at once for the Index/Date columns
$ws.Range("A10").Value = 3
$ws.Range("B10").Value = "31May2020  2330"

$ws.Range("A11").Value = 4
$ws.Range("B11").Value = "31May2020  2345"

$ws.Range("A12").Value = 5
$ws.Range("B12").Value = "01Jun2020  0000"

$ws.Range("A13").Value = 6
$ws.Range("B13").Value = "01Jun2020  0015"

$ws.Range("A14").Value = 7
$ws.Range("B14").Value = "01Jun2020  0030"

$ws.Range("A15").Value = 8
$ws.Range("B15").Value = "01Jun2020  0045"

$ws.Range("A16").Value = 9
$ws.Range("B16").Value = "01Jun2020  0100"

$ws.Range("A17").Value = 10
$ws.Range("B17").Value = "01Jun2020  0115"

# Value column for rows 10-17: a single relative formula assigned to the whole
# range at once creates one shared formula group (matches t="shared" si="0").
$ws.Range("C10:C17").Formula = "=C9+1"

# Update selection to match target (row 8 selected entirely)
$ws.Range("A8:XFD8").Select()
